$d = $word.ActiveDocument

# --- Merge split runs "3"+". frmThongKe" -> "3. frmThongKe" (and similar) ---
# Word's Find/Execute with Replace=wdReplaceAll (2) collapses the matched
# range (which spans the two adjacent runs) into a single run carrying the
# found run's formatting - exactly mirroring the diff's run-merge.
$mergeParagraphIndexes = @(33, 34, 44, 45)
foreach ($idx in $mergeParagraphIndexes) {
    $p = $d.Paragraphs.Item($idx)
    $txt = $p.Range.Text
    $txt = $txt.TrimEnd([char]13, [char]7)
    $rng = $p.Range
    $rng.Find.Execute($txt, $false, $false, $false, $false, $false, $true, 1, $false, $txt, 2)
}

# --- Insert two new paragraphs after the first empty (720-indent) paragraph ---
# Before: ...Ae co gi... | empty(720) | empty(720) | tab-only...
# After:  ...Ae co gi... | empty(720) | empty(720, NEW) | "Xuaan Dduc commit"(720, NEW) | empty(720) | tab-only...
$firstEmpty = $d.Paragraphs.Item(51)
$insertPos = $firstEmpty.Range.End
$twoBreaks = "$([char]13)$([char]13)"
$d.Range($insertPos, $insertPos).Text = $twoBreaks

$target = $d.Paragraphs.Item(53)
$targetRange = $target.Range
$targetRange.Text = "Xuaan Dduc commit"
$targetRange = $target.Range
$targetRange.Font.Size = 14
$targetRange.Font.SizeBi = 14
